$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the existing hyperlink on C5 (row 5 no longer links out) before
# changing the cell's text, but keep the cell's "hyperlink" font style.
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Address($false, $false) -eq "C5") {
        $hl.Delete()
    }
}

# Update the data row (row 5) to the new "gw" environment values.
# B5 keeps its "quote-prefixed text" look (same as the other rows in
# this column) - a leading apostrophe tells Excel to keep storing it
# as explicit text without altering the cell's format.
$ws.Range("B5").Value = "'ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("C5").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("E5").Value = "gw"
$ws.Range("F5").Value = 8684079401
$ws.Range("G5").Value = 2302
$ws.Range("H5").Value = "Mattioli"
$ws.Range("AA5").Value = "RGM104"
$ws.Range("AB5").Value = "ABC12SRGM104"
$ws.Range("AC5").Value = "ZAZ123SRGM104"

# G5/H5 drop the custom (Arial) formatting that used to mark them and
# fall back to the sheet's normal style, matching the other "gw" rows.
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Style = "Normal"

# Move the current selection to the full row 5 (no frozen/scrolled
# top-left column anymore).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows("5").Select()
